$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values following a data repull / recalculation
$ws.Range("F2").Value = -3
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 5
$ws.Range("F7").Value = 8
